$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0283232344
$ws.Range("C2").Value = 0.0006822424
$ws.Range("D2").Value = 0.0004153632
$ws.Range("E2").Value = 0.0009877962000000001
$ws.Range("F2").Value = 0.0007164220000000001
$ws.Range("G2").Value = 0.000424609

$ws.Range("B3").Value = 0.2121381092
$ws.Range("C3").Value = 0.0032709344
$ws.Range("D3").Value = 0.0008710728
$ws.Range("E3").Value = 0.0019942228
$ws.Range("F3").Value = 0.0027361546
$ws.Range("G3").Value = 0.0009780699999999999

$ws.Range("B4").Value = 1.6764079678
$ws.Range("C4").Value = 0.013376783
$ws.Range("D4").Value = 0.0018909846
$ws.Range("E4").Value = 0.0074744184
$ws.Range("F4").Value = 0.0118008174
$ws.Range("G4").Value = 0.0024347038

$ws.Range("B5").Value = 14.0336218528
$ws.Range("C5").Value = 0.0517733154
$ws.Range("D5").Value = 0.0041477716
$ws.Range("E5").Value = 0.0372150808
$ws.Range("F5").Value = 0.05080168500000001
$ws.Range("G5").Value = 0.0060820926

$ws.Range("B6").Value = 116.941428585
$ws.Range("C6").Value = 0.1998502272
$ws.Range("D6").Value = 0.0094201748
$ws.Range("E6").Value = 0.109190961
$ws.Range("F6").Value = 0.1903564716
$ws.Range("G6").Value = 0.0131280356

$ws.Range("B7").Value = 966.9290217122
$ws.Range("C7").Value = 0.8066466206
$ws.Range("D7").Value = 0.0207859002
$ws.Range("E7").Value = 0.3700318532
$ws.Range("F7").Value = 0.7497630564000001
$ws.Range("G7").Value = 0.028741423

$ws.Range("C8").Value = 3.2098342538
$ws.Range("D8").Value = 0.0450456594
$ws.Range("E8").Value = 1.5996798728
$ws.Range("F8").Value = 2.9972841358
$ws.Range("G8").Value = 0.0644687556
